$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BioSample")
$ws.Columns("AA").Delete()
